$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended after row 233 - update through 2021-04-26 (dates 2021-04-22..2021-04-26)
$newRows = @(
    @{ Row = 234; A = 44308; B = 1; C = 3; D = 327.5109170305677 },
    @{ Row = 235; A = 44309; B = 1; C = 3; D = 327.5109170305677 },
    @{ Row = 236; A = 44310; B = 1; C = 4; D = 436.6812227074236 },
    @{ Row = 237; A = 44311; B = 0; C = 4; D = 436.6812227074236 },
    @{ Row = 238; A = 44312; B = 0; C = 4; D = 436.6812227074236 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A keeps the same date-number-format/bold/centered style as the
    # rows above it; copy that formatting down before writing the new value.
    $ws.Range("A233").Copy()
    $ws.Range("A$rowNum").PasteSpecial(-4122)

    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
}
